$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the readonly text box in info sheet tab: the "Ingredients" and
# "Used By & Best Before Date" columns (F and G) had their header/value
# pairs mismatched. Swap columns F and G (both header row 1 and value row 2)
# so headers line up with their correct values again.

$f1 = $ws.Range("F1").Text
$g1 = $ws.Range("G1").Text
$f2 = $ws.Range("F2").Text
$g2 = $ws.Range("G2").Text

$ws.Range("F1").Value = $g1
$ws.Range("G1").Value = $f1
$ws.Range("F2").Value = $g2
$ws.Range("G2").Value = $f2
